$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.623.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.098.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.099.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.490"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.596.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.419.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.087.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "495.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.701"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.73%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.88%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.93%  "

$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0416"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "448.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.066.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("E44").Value = "  +8.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.07%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("E48").Value = "  +1.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0533"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.23%  "
